$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.81 = 40432.19 pesos`n✅ 40432.19 pesos = 9.79 = 935.16 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the tasas figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 101.899
$wsTasas.Range("O10").Value = 4120
$wsTasas.Range("N12").Value = 4129
$wsTasas.Range("O12").Value = 95.5
